$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values
# (e.g. "289.85", "-3.75%") rather than numeric/percentage values, so the
# target cells must be forced to Text format before assignment. Otherwise
# Excel would auto-convert the strings to numbers and mangle formatting
# (trailing zeros, percent scaling, etc).
$cellEdits = @(
    @{ Cell = "D2";  Value = "289.85" }
    @{ Cell = "E2";  Value = "-3.75%" }
    @{ Cell = "D3";  Value = "30.90" }
    @{ Cell = "E3";  Value = "-3.90%" }
    @{ Cell = "D4";  Value = "4.868" }
    @{ Cell = "E4";  Value = "-2.35%" }
    @{ Cell = "D5";  Value = "0.07153" }
    @{ Cell = "E5";  Value = "-9.39%" }
    @{ Cell = "D6";  Value = "1.832" }
    @{ Cell = "E6";  Value = "-12.62%" }
    @{ Cell = "D7";  Value = "7.644" }
    @{ Cell = "E7";  Value = "-1.95%" }
    @{ Cell = "D8";  Value = "3.770" }
    @{ Cell = "E8";  Value = "-1.55%" }
    @{ Cell = "D9";  Value = "0.8935" }
    @{ Cell = "E9";  Value = "-3.74%" }
    @{ Cell = "D10"; Value = "0.1643" }
    @{ Cell = "E10"; Value = "-5.87%" }
    @{ Cell = "D11"; Value = "0.07533" }
    @{ Cell = "E11"; Value = "-5.32%" }
    @{ Cell = "D12"; Value = "0.08103" }
    @{ Cell = "E12"; Value = "-6.51%" }
    @{ Cell = "D13"; Value = "0.02984" }
    @{ Cell = "E13"; Value = "-3.96%" }
    @{ Cell = "E14"; Value = "-0.17%" }
    @{ Cell = "D15"; Value = "0.001490" }
    @{ Cell = "E15"; Value = "-1.84%" }
    @{ Cell = "D16"; Value = "0.005843" }
    @{ Cell = "E16"; Value = "-2.23%" }
    @{ Cell = "D18"; Value = "3.468" }
    @{ Cell = "E18"; Value = "0.21%" }
    @{ Cell = "D19"; Value = "2.106" }
    @{ Cell = "E19"; Value = "-7.47%" }
    @{ Cell = "D20"; Value = "0.3278" }
    @{ Cell = "E20"; Value = "-0.29%" }
    @{ Cell = "E21"; Value = "-1.30%" }
    @{ Cell = "D22"; Value = "4.265" }
    @{ Cell = "E22"; Value = "-0.35%" }
    @{ Cell = "D23"; Value = "0.1999" }
    @{ Cell = "E23"; Value = "11.53%" }
    @{ Cell = "D24"; Value = "0.04469" }
    @{ Cell = "E24"; Value = "-3.03%" }
    @{ Cell = "D25"; Value = "0.001211" }
    @{ Cell = "E25"; Value = "-2.06%" }
    @{ Cell = "D26"; Value = "0.004663" }
    @{ Cell = "E26"; Value = "4.87%" }
    @{ Cell = "D27"; Value = "0.0001250" }
    @{ Cell = "E27"; Value = "-0.14%" }
    @{ Cell = "D39"; Value = "0.01638" }
    @{ Cell = "E39"; Value = "-4.59%" }
    @{ Cell = "D40"; Value = "0.04358" }
    @{ Cell = "E40"; Value = "-8.52%" }
    @{ Cell = "D41"; Value = "0.007383" }
    @{ Cell = "E41"; Value = "-0.78%" }
    @{ Cell = "E42"; Value = "-3.92%" }
    @{ Cell = "D43"; Value = "0.002003" }
    @{ Cell = "E43"; Value = "-11.89%" }
    @{ Cell = "D44"; Value = "0.01016" }
    @{ Cell = "E44"; Value = "-10.18%" }
    @{ Cell = "D45"; Value = "0.00005844" }
    @{ Cell = "E45"; Value = "-2.35%" }
    @{ Cell = "D46"; Value = "0.00000000750" }
    @{ Cell = "E46"; Value = "-0.13%" }
    @{ Cell = "D47"; Value = "2.210" }
    @{ Cell = "E47"; Value = "168.38%" }
    @{ Cell = "E48"; Value = "-11.48%" }
    @{ Cell = "D49"; Value = "0.00002100" }
    @{ Cell = "E49"; Value = "-0.13%" }
    @{ Cell = "D50"; Value = "0.0002000" }
    @{ Cell = "E50"; Value = "-0.13%" }
)

foreach ($edit in $cellEdits) {
    $range = $ws.Range($edit.Cell)
    $range.NumberFormat = "@"
    $range.Value = $edit.Value
}
